$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.202020202020202
$ws.Range("C2").Value = 0.5555555555555556
$ws.Range("J2").Value = 0.04040404040404041
$ws.Range("P2").Value = 0.1414141414141414
$ws.Range("S2").Value = 0.06060606060606061
$ws.Range("C3").Value = 0.03333333333333333
$ws.Range("J3").Value = 0.01666666666666667
$ws.Range("P3").Value = 0.7833333333333333
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.04950495049504951
$ws.Range("D6").Value = 0.009900990099009901
$ws.Range("F6").Value = 0.09900990099009901
$ws.Range("J6").Value = 0.2574257425742574
$ws.Range("O6").Value = 0.06930693069306931
$ws.Range("Q6").Value = 0.1188118811881188
$ws.Range("R6").Value = 0.0396039603960396
$ws.Range("S6").Value = 0.3564356435643564
$ws.Range("B7").Value = 0.1551724137931035
$ws.Range("D7").Value = 0.05172413793103448
$ws.Range("F7").Value = 0.06896551724137931
$ws.Range("J7").Value = 0.1379310344827586
$ws.Range("O7").Value = 0.01724137931034483
$ws.Range("Q7").Value = 0.1896551724137931
$ws.Range("R7").Value = 0.01724137931034483
$ws.Range("S7").Value = 0.3620689655172414
$ws.Range("B8").Value = 0.1375
$ws.Range("D8").Value = 0.0125
$ws.Range("F8").Value = 0.05
$ws.Range("J8").Value = 0.1
$ws.Range("O8").Value = 0.0625
$ws.Range("Q8").Value = 0.125
$ws.Range("R8").Value = 0.0625
$ws.Range("S8").Value = 0.45
$ws.Range("B9").Value = 0.04761904761904762
$ws.Range("D9").Value = 0.03174603174603174
$ws.Range("F9").Value = 0.1111111111111111
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("O9").Value = 0.07936507936507936
$ws.Range("Q9").Value = 0.2063492063492063
$ws.Range("R9").Value = 0.06349206349206349
$ws.Range("S9").Value = 0.3174603174603174
$ws.Range("B10").Value = 0.1297297297297297
$ws.Range("D10").Value = 0.03783783783783784
$ws.Range("F10").Value = 0.0945945945945946
$ws.Range("J10").Value = 0.1189189189189189
$ws.Range("O10").Value = 0.02972972972972973
$ws.Range("Q10").Value = 0.2216216216216216
$ws.Range("R10").Value = 0.05945945945945946
$ws.Range("S10").Value = 0.3081081081081081
$ws.Range("G11").Value = 0.1739130434782609
$ws.Range("J11").Value = 0.108695652173913
$ws.Range("K11").Value = 0.2173913043478261
$ws.Range("L11").Value = 0.5
$ws.Range("G12").Value = 0.723404255319149
$ws.Range("J12").Value = 0.1914893617021277
$ws.Range("K12").Value = 0.02127659574468085
$ws.Range("S12").Value = 0.06382978723404255
$ws.Range("G13").Value = 0.8333333333333334
$ws.Range("J13").Value = 0.1666666666666667
$ws.Range("F15").Value = 0.02352941176470588
$ws.Range("H15").Value = 0.08235294117647059
$ws.Range("I15").Value = 0.04705882352941176
$ws.Range("J15").Value = 0.3411764705882353
$ws.Range("K15").Value = 0.08235294117647059
$ws.Range("M15").Value = 0.01176470588235294
$ws.Range("O15").Value = 0.1176470588235294
$ws.Range("S15").Value = 0.2941176470588235
$ws.Range("F16").Value = 0.04166666666666666
$ws.Range("H16").Value = 0.04166666666666666
$ws.Range("I16").Value = 0.06944444444444445
$ws.Range("J16").Value = 0.4861111111111111
$ws.Range("K16").Value = 0.1666666666666667
$ws.Range("M16").Value = 0.04166666666666666
$ws.Range("O16").Value = 0.04166666666666666
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.03937007874015748
$ws.Range("H17").Value = 0.1338582677165354
$ws.Range("I17").Value = 0.1102362204724409
$ws.Range("J17").Value = 0.4330708661417323
$ws.Range("K17").Value = 0.1102362204724409
$ws.Range("O17").Value = 0.07086614173228346
$ws.Range("S17").Value = 0.1023622047244094
$ws.Range("F18").Value = 0.02857142857142857
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.1714285714285714
$ws.Range("J18").Value = 0.2857142857142857
$ws.Range("K18").Value = 0.08571428571428572
$ws.Range("M18").Value = 0.02857142857142857
$ws.Range("O18").Value = 0.1142857142857143
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.06116207951070336
$ws.Range("H19").Value = 0.1406727828746177
$ws.Range("I19").Value = 0.1009174311926606
$ws.Range("J19").Value = 0.3883792048929663
$ws.Range("K19").Value = 0.1100917431192661
$ws.Range("M19").Value = 0.03363914373088685
$ws.Range("O19").Value = 0.09174311926605505
$ws.Range("S19").Value = 0.07339449541284404
